$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column J describes "provincia". It was previously modeled as a
# dimension referencing "sdmx-dimension:refArea" with a "URI-Provincia"
# type; it is now re-curated as a plain measure ("iaest-measure:provincia")
# typed as "xsd:int", matching how the other measure columns are modeled.

$ws.Range("J2").Value = "iaest-measure:provincia"
$ws.Range("J3").Value = "medida"
$ws.Range("J4").Value = "xsd:int"
